$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1) Duplicate "Betão F" into a new sheet "Carlão", placed right before "Betão F"
#    (so the tab order becomes: Menu, Aplicação, Carlão, Betão F, Graficos, dados, ...)
$betaoF = $wb.Worksheets.Item("Betão F")
$betaoF.Copy($betaoF)
$carlao = $wb.Worksheets.Item("Betão F (2)")
$carlao.Name = "Carlão"

# 2) Give the new "Carlão" sheet its own sheet-scoped "Betao_sal" defined name,
#    mirroring the workbook-scoped one that still points at "Betão F".
$carlao.Names.Add("Betao_sal", "=Carlão!`$I`$8:`$I`$37")

# 3) Restore the view state on each worksheet that the author left behind.
$carlao.Activate()
$carlao.Range("L27").Select()

$betaoF.Activate()
$betaoF.Range("K32").Select()

# 4) Drop the unused blank placeholder sheets (Plan7 .. Plan12).
foreach ($name in @("Plan7", "Plan8", "Plan9", "Plan10", "Plan11", "Plan12")) {
    $wb.Worksheets.Item($name).Delete()
}

# 5) Leave "dados" as the active/visible tab, matching the final saved state.
$dados = $wb.Worksheets.Item("dados")
$dados.Activate()
$dados.Range("E11").Select()
